$d = $word.ActiveDocument
$tbl = $d.Tables.Item(2)

# Row 2, Cell 4 (State) - simple replace Incomplete -> Complete
$cell2 = $tbl.Rows.Item(2).Cells.Item(4)
$rng2 = $cell2.Range
$rng2.Find.Execute("Incomplete", $false, $false, $false, $false, $false, $true, 1, $false, "Complete", 2)

# Row 3, Cell 4 (State) - simple replace Incomplete -> Complete
$cell3 = $tbl.Rows.Item(3).Cells.Item(4)
$rng3 = $cell3.Range
$rng3.Find.Execute("Incomplete", $false, $false, $false, $false, $false, $true, 1, $false, "Complete", 2)

Write-Output "Done rows 2,3"

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $row = $tbl.Rows.Item($r)
    $line = ""
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cellText = $row.Cells.Item($c).Range.Text
        $line = $line + "[" + $cellText + "]"
    }
    Write-Output "Row $r : $line"
}
